$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 7 holds the "06" entry: B7 = "Ringgold", C7 = "Ringgold identifier..."
# Fix the misspelling "Ringgold" -> "Ringold" in both cells.
$ws.Range("B7").Value = "Ringold"
$ws.Range("C7").Value = "Ringold identifier for organisations in the publishing industry supply chain"
